# Append two new price-tracking rows to the sheet (row 64 and row 65),
# matching the existing "data / hora / preco / site / cor" layout, and
# carrying forward the date-time number format used by the rest of
# column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFormat = "YYYY-MM-DD HH:MM:SS"

# Row 64: amazon, preto
$ws.Range("A64").Value = 45211
$ws.Range("A64").NumberFormat = $dateFormat
$ws.Range("B64").Value = "16:11"
$ws.Range("C64").Value = 2552
$ws.Range("D64").Value = "amazon"
$ws.Range("E64").Value = "preto"

# Row 65: mercado livre, preto
$ws.Range("A65").Value = 45211
$ws.Range("A65").NumberFormat = $dateFormat
$ws.Range("B65").Value = "16:12"
$ws.Range("C65").Value = 2563
$ws.Range("D65").Value = "mercado livre"
$ws.Range("E65").Value = "preto"
